$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.611.05"
$ws.Range("E2").Value = "  +1.47%  "

$ws.Range("D3").Value = "3.496.05"
$ws.Range("E3").Value = "  +0.92%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'591.70"
$ws.Range("E5").Value = "  +1.56%  "

$ws.Range("D6").Value = "'168.62"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = "  +5.16%  "

$ws.Range("D9").Value = "'0.129"
$ws.Range("E9").Value = "  +4.54%  "

$ws.Range("E10").Value = "  +0.37%  "

$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").Value = "4.099.88"
$ws.Range("E12").Value = "  +1.01%  "

$ws.Range("E13").Value = "  -0.23%  "

$ws.Range("D14").Value = "'28.13"
$ws.Range("E14").Value = "  +2.26%  "

$ws.Range("D15").Value = "'0.0000179"
$ws.Range("E15").Value = "  +1.33%  "

$ws.Range("D16").Value = "66.660.80"
$ws.Range("E16").Value = "  +1.69%  "

$ws.Range("D17").Value = "3.495.09"
$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("D18").Value = "'6.30"
$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("D19").Value = "'14.05"
$ws.Range("E19").Value = "  +2.08%  "

$ws.Range("D20").Value = "'393.68"
$ws.Range("E20").Value = "  +2.62%  "

$ws.Range("D21").Value = "'7.92"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").Value = "'73.01"
$ws.Range("E22").Value = "  +1.67%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("E24").Value = "  +2.67%  "

$ws.Range("E25").Value = "  +1.22%  "

$ws.Range("D26").Value = "'10.14"
$ws.Range("E26").Value = "  +3.42%  "

$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").Value = "'6.35"
$ws.Range("E29").Value = "  +1.89%  "

$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("D31").Value = "'2.06"
$ws.Range("E31").Value = "  +1.51%  "

$ws.Range("D32").Value = "'23.76"
$ws.Range("E32").Value = "  +1.94%  "

$ws.Range("E33").Value = "  +0.54%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'1.60"
$ws.Range("E34").Value = "  +5.03%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'162.66"
$ws.Range("E35").Value = "  +1.56%  "

$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").Value = "'0.896"
$ws.Range("E36").Value = "  +0.63%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'1.91"
$ws.Range("E37").Value = "  +2.24%  "

$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'6.80"
$ws.Range("E38").Value = "  +3.01%  "

$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'4.65"
$ws.Range("E39").Value = "  +4.15%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'26.54"
$ws.Range("E40").Value = "  +1.70%  "

$ws.Range("D41").Value = "'0.0740"
$ws.Range("E41").Value = "  +0.51%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'26.78"
$ws.Range("E42").Value = "  +0.21%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.772.23"
$ws.Range("E43").Value = "  -1.03%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'42.94"
$ws.Range("E44").Value = "  -0.25%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.55"
$ws.Range("E45").Value = "  +2.39%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0310"
$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "'342.16"
$ws.Range("E47").Value = "  +1.63%  "

$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.09"
$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "'33.95"
$ws.Range("E49").Value = "  +4.38%  "

$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "'0.854"
$ws.Range("E50").Value = "  +2.86%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'6.49"
$ws.Range("E51").Value = "  +1.55%  "

